$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: update title and link
$ws.Range("D9").Value = "MBA AI/BigData 과정 입학시험 문제 예시"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/mba-ai-admission-exam-sample/#utm_source=rss&utm_medium=rss&utm_campaign=mba-ai-admission-exam-sample"

# Row 29: update title
$ws.Range("D29").Value = "[만화] 인턴일기 72~80"

# Row 50: update title and link
$ws.Range("D50").Value = "RND 100선 투표 [ 추첨 200명 ]"
$ws.Range("E50").Value = "http://incredible.egloos.com/7572335"

# Row 51: update title and link
$ws.Range("D51").Value = "[MySQL] 스토어드 프로시저로 데이터 insert 하기"
$ws.Range("E51").Value = "https://bskyvision.com/entry/MySQL-%EC%8A%A4%ED%86%A0%EC%96%B4%EB%93%9C-%ED%94%84%EB%A1%9C%EC%8B%9C%EC%A0%80%EB%A1%9C-%EB%8D%B0%EC%9D%B4%ED%84%B0-insert-%ED%95%98%EA%B8%B0"
